# Apply the FHIR StructureDefinition metadata refresh:
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to the new publication timestamp
#  - set Publisher to "Alvearie Team"
#  - replace the (duplicated) "Contact" row with a single "Jurisdiction" row
#  - update the Extension element's Short/Definition text on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Remove the second, duplicate "Contact" row (row 11); this shifts every
# subsequent row up by one, turning the old A1:B21 range into A1:B20.
$meta.Rows.Item(11).Delete()

# Property/value updates on the Metadata sheet.
$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the top-level "Extension" element; give it the real short
# description / definition instead of the generic Extension boilerplate.
$elements.Range("K2").Value = "Employee Sub Business Unit"
$elements.Range("L2").Value = "Code of the sub business unit of the employee"
